# Daily attendance processing - 2026-01-01 17:57:37
# For every row in column G ("Recorded By"), if the comma-separated list of
# recorders contains "System" but "System" is not already the first entry,
# reverse the order of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Length -gt 0) {
        $parts = $val -split ', '
        if ($parts.Count -gt 1 -and ($parts -contains 'System') -and ($parts[0] -ne 'System')) {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(', ', $reversed)
        }
    }
}
